$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.12831060389814297
$ws.Range("A2").Value = -0.0059999999887949684
$ws.Range("A3").Value = 0.0021970331440908097
$ws.Range("A4").Value = -0.0079999999829141188
$ws.Range("A5").Value = -0.0029999999914291919
$ws.Range("A6").Value = -0.0019999999916802125
$ws.Range("A7").Value = -0.0099999999767317327
$ws.Range("A8").Value = -0.0099999999759874392
$ws.Range("A9").Value = -0.0019999999898292486
$ws.Range("A10").Value = 0.0087802852090987926
$ws.Range("A11").Value = -0.0029999999871668237
$ws.Range("A12").Value = -0.0034999999859830488
$ws.Range("A13").Value = -0.0034999999847427077
$ws.Range("A14").Value = -0.0079999999761790619
$ws.Range("A15").Value = -0.00099999998877464691
$ws.Range("A16").Value = -0.0019999999867450491
$ws.Range("A17").Value = -0.0019999999864470652
$ws.Range("A18").Value = -0.0039999999827058375
$ws.Range("A19").Value = -0.0039999999923896468
$ws.Range("A20").Value = -0.0039999999917981199
$ws.Range("A21").Value = -0.0039999999917261775
$ws.Range("A22").Value = -0.003999999991654235
$ws.Range("A23").Value = -0.0049999999884873247
$ws.Range("A24").Value = -0.019999999959854797
$ws.Range("A25").Value = -0.019999999959311232
$ws.Range("A26").Value = -0.0024999999907233672
$ws.Range("A27").Value = -0.0024999999907127091
$ws.Range("A28").Value = -0.001999999991543433
$ws.Range("A29").Value = 0.010066299353738373
$ws.Range("A30").Value = -0.059999999885437916
$ws.Range("A31").Value = 0.06673091813475196
$ws.Range("A32").Value = -0.0099999999785147509
$ws.Range("A33").Value = -0.0039999999894462235
